$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells in column D that look like plain numbers need to be forced back
# to Text (matching the original inlineStr cell type) before assigning,
# otherwise Excel auto-converts them to numeric values (e.g. dropping
# trailing zeros: "87.20" -> 87.2).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '39.229.24'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').Value = '2.201.10'
$ws.Range('E3').Value = '  -5.45%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '295.11'
$ws.Range('E5').Value = '  -3.83%  '
$ws.Range('D6').Value = '81.03'
$ws.Range('E6').Value = '  -4.23%  '
$ws.Range('E7').Value = '  -3.86%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('E10').Value = '  -5.53%  '
$ws.Range('D11').Value = '29.07'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').Value = '46.83'
$ws.Range('E12').Value = '  -11.19%  '
$ws.Range('E13').Value = '  -2.65%  '
$ws.Range('D14').Value = '2.532.96'
$ws.Range('E14').Value = '  -5.72%  '
$ws.Range('D15').Value = '6.21'
$ws.Range('E15').Value = '  -3.04%  '
$ws.Range('E16').Value = '  -4.96%  '
$ws.Range('D17').Value = '2.196.33'
$ws.Range('E17').Value = '  -5.69%  '
$ws.Range('E18').Value = '  -5.46%  '
$ws.Range('D19').Value = '39.147.85'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('E20').Value = '  -3.47%  '
$ws.Range('E21').Value = '  -6.11%  '
$ws.Range('D22').Value = '64.47'
$ws.Range('E22').Value = '  -4.60%  '
$ws.Range('D23').Value = '10.22'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').Value = '226.39'
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  -5.90%  '
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.25'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '22.50'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('D30').Value = '9.03'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').Value = '149.39'
$ws.Range('E31').Value = '  -1.71%  '
$ws.Range('D32').Value = '31.47'
$ws.Range('E32').Value = '  -11.20%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').Value = '4.78'
$ws.Range('E34').Value = '  -6.25%  '
$ws.Range('E35').Value = '  -4.33%  '
$ws.Range('D36').Value = '0.0692'
$ws.Range('E36').Value = '  -4.08%  '
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('D38').Value = '15.28'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('D39').Value = '0.0959'
$ws.Range('E39').Value = '  -3.65%  '
$ws.Range('E40').Value = '  -4.63%  '
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('D42').Value = '3.60'
$ws.Range('E42').Value = '  -5.65%  '
$ws.Range('D43').Value = '1.897.51'
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('E44').Value = '  -10.29%  '
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').Value = '16.04'
$ws.Range('E47').Value = '  -8.49%  '
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('D49').Value = '71.14'
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').Value = '2.402.73'
$ws.Range('E50').Value = '  -6.05%  '
$ws.Range('D51').Value = '87.20'
$ws.Range('E51').Value = '  -5.93%  '
